$wb = $excel.ActiveWorkbook

# --- Sheet "12L": just move the selection (no content changes) ---
$ws12L = $wb.Worksheets.Item("12L")
$ws12L.Range("C4:C18").Select()

# --- Sheet "Linh Tinh": just move the selection (no content changes) ---
$wsLinhTinh = $wb.Worksheets.Item("Linh Tinh")
$wsLinhTinh.Range("C13:C18").Select()

# --- Sheet "RVC": remove "Pham Minh Triet" row, shift names up, add "Nghia Le" at the end ---
$wsRVC = $wb.Worksheets.Item("RVC")

$wsRVC.Range("C7").Value = "Lê Thị Bích Phượng"
$wsRVC.Range("C8").Value = "Nguyễn Ngọc Thạch"
$wsRVC.Range("C9").Value = "Đặng Khánh Toàn"
$wsRVC.Range("C10").Value = "Phạm Hoàng Mai"
$wsRVC.Range("C11").Value = "Nguyễn Hoàng Kiên"
$wsRVC.Range("C12").Value = "Khoa Nguyễn"
$wsRVC.Range("C13").Value = "Lê Văn Quân"
$wsRVC.Range("C14").Value = "Phó Kiến Huy"
$wsRVC.Range("C15").Value = "Giang Châu"
$wsRVC.Range("C16").Value = "Đặng Trần Công Lý"
$wsRVC.Range("C17").Value = "Hiếu Nguyễn"
$wsRVC.Range("C18").Value = "Thành Phạm"
$wsRVC.Range("C19").Value = "Nghĩa Lê"

# Mark the first guest as "No" (not attending)
$wsRVC.Range("D3").Value = "No"

# Update the Yes/No validation list to also allow "1" and "2"
$wsRVC.Range("D3:D19").Validation.Modify(3, 1, 1, """1, 2, No""")

# Leave the final selection on RVC, matching the saved workbook state
$wsRVC.Range("D5").Select()
